$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.529.19"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").Value = "3.429.61"
$ws.Range("E3").Value = "  +2.42%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'407.44"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").Value = "'130.80"
$ws.Range("E6").Value = "  +3.18%  "

$ws.Range("D7").Value = "'0.598"
$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.695"
$ws.Range("E9").Value = "  +5.44%  "

$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  +19.26%  "

$ws.Range("D11").Value = "'42.20"
$ws.Range("E11").Value = "  +2.91%  "

$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "'8.48"
$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("D14").Value = "'19.88"
$ws.Range("E14").Value = "  +2.96%  "

$ws.Range("D15").Value = "3.434.96"
$ws.Range("E15").Value = "  +2.38%  "

$ws.Range("D16").Value = "62.592.72"
$ws.Range("E16").Value = "  +1.98%  "

$ws.Range("E17").Value = "  +2.42%  "

$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("D19").Value = "'0.0000164"
$ws.Range("E19").Value = "  +29.03%  "

$ws.Range("D20").Value = "'3.19"
$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").Value = "'84.52"
$ws.Range("E21").Value = "  +5.56%  "

$ws.Range("D22").Value = "'314.55"
$ws.Range("E22").Value = "  +5.01%  "

$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").Value = "'3.18"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").Value = "'4.74"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").Value = "'29.84"
$ws.Range("E26").Value = "  +3.03%  "

$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("E28").Value = "  +5.17%  "

$ws.Range("D29").Value = "'2.73"
$ws.Range("E29").Value = "  +10.01%  "

$ws.Range("D30").Value = "'44.52"
$ws.Range("E30").Value = "  +8.87%  "

$ws.Range("D31").Value = "'0.174"
$ws.Range("E31").Value = "  +1.62%  "

$ws.Range("D32").Value = "'0.115"
$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("D33").Value = "'11.44"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").Value = "'0.0485"
$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("D36").Value = "'51.81"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "'2.97"
$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("D39").Value = "'0.323"
$ws.Range("E39").Value = "  +15.67%  "

$ws.Range("D40").Value = "'3.34"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("D41").Value = "'143.63"
$ws.Range("E41").Value = "  +5.21%  "

$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("D43").Value = "'1.99"
$ws.Range("E43").Value = "  +1.24%  "

$ws.Range("D44").Value = "'16.95"
$ws.Range("E44").Value = "  +1.60%  "

$ws.Range("D45").Value = "'3.93"
$ws.Range("E45").Value = "  +1.47%  "

$ws.Range("D46").Value = "'2.23"
$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("D47").Value = "'21.39"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").Value = "2.112.11"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").Value = "'2.00"
$ws.Range("E49").Value = "  +8.42%  "

$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("E51").Value = "  +30.21%  "
